# إضافة حدث جديد في Card20 by admin at 2025-12-08 09:13:39
# Fill the previously-blank filler cells on the last existing event row (23)
# with "nan" (matching the convention used by every other event row), then
# append a brand-new service-event row (24) for the card.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card20")

# --- Row 23: the blank placeholder cells get the "nan" filler text ---
$nanCols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "M")
foreach ($col in $nanCols) {
    $ws.Range($col + "23").Value = "nan"
}

# --- Row 24: new service event ---
# Column A ("card") mirrors every other row's text "20" label. A plain
# Value assignment would be auto-coerced to a number by Excel, so force
# text entry via a temporary Text number format, then drop the formatting
# again so no stray style lingers on the cell.
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "20"
$ws.Range("A24").ClearFormats()

# B24:K24 stay empty, exactly like row 23 did before this edit.

$ws.Range("L24").Value = "14\8\2025"
$ws.Range("M24").Value = "766 t"
$ws.Range("N24").Value = "تم تغيير زيت الجيربوكس"
$ws.Range("O24").Value = "تيم العمل"
